$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data of rows 40 and 41 (match order / scrape order changed) ---
$row40 = @("Balzan", 0, "Birkirkara", 3, 2.8, "28/10/2023 09:43", 3.6, "29/10/2023 14:51", 3.03, "28/10/2023 09:43", 2.56, "29/10/2023 14:51", 2.34, "28/10/2023 09:43", 2.49, "29/10/2023 14:51", "https://www.betexplorer.com/football/malta/premier-league/balzan-fc-birkirkara/WnBN3OYq/")
$row41 = @("Gudja", 1, "Floriana", 1, 4.82, "28/10/2023 09:43", 9.960000000000001, "29/10/2023 14:43", 3.62, "28/10/2023 09:43", 5.01, "29/10/2023 14:43", 1.58, "28/10/2023 09:43", 1.3, "29/10/2023 10:32", "https://www.betexplorer.com/football/malta/premier-league/gudja-floriana/UwPe84QS/")

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "40").Value = $row41[$i]
    $ws.Range($cols[$i] + "41").Value = $row40[$i]
}

# --- Append two new rows (48, 49) for matches added by the scraper ---
# Copy row 47's formatting down so the "Indice" (col A) and "data_partida" (col E)
# styles match the rest of the table, then overwrite the values.
$ws.Range("A47:V47").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "malta"
$ws.Range("C48").Value = "premier-league"
$ws.Range("D48").Value = "2023-2024"
$ws.Range("E48").Value = 45235.6875
$ws.Range("F48").Value = "Balzan"
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = "Naxxar"
$ws.Range("I48").Value = 2
$ws.Range("J48").Value = 1.38
$ws.Range("K48").Value = "04/11/2023 04:43"
$ws.Range("L48").Value = 1.53
$ws.Range("M48").Value = "05/11/2023 12:57"
$ws.Range("N48").Value = 4.46
$ws.Range("O48").Value = "04/11/2023 04:43"
$ws.Range("P48").Value = 4.16
$ws.Range("Q48").Value = "05/11/2023 14:31"
$ws.Range("R48").Value = 5.95
$ws.Range("S48").Value = "04/11/2023 04:43"
$ws.Range("T48").Value = 5.47
$ws.Range("U48").Value = "05/11/2023 12:57"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/malta/premier-league/balzan-fc-naxxar-lions/8r3YLNfS/"

$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "malta"
$ws.Range("C49").Value = "premier-league"
$ws.Range("D49").Value = "2023-2024"
$ws.Range("E49").Value = 45235.6875
$ws.Range("F49").Value = "Floriana"
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = "Sirens"
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1.3
$ws.Range("K49").Value = "04/11/2023 04:43"
$ws.Range("L49").Value = 1.38
$ws.Range("M49").Value = "05/11/2023 16:21"
$ws.Range("N49").Value = 4.68
$ws.Range("O49").Value = "04/11/2023 04:43"
$ws.Range("P49").Value = 3.9
$ws.Range("Q49").Value = "05/11/2023 16:24"
$ws.Range("R49").Value = 7.42
$ws.Range("S49").Value = "04/11/2023 04:43"
$ws.Range("T49").Value = 11.61
$ws.Range("U49").Value = "05/11/2023 16:24"
$ws.Range("V49").Value = "https://www.betexplorer.com/football/malta/premier-league/floriana-sirens/Ui4UM3uM/"
